# Transect Setup Data Sheet - add Panoche site transect rows (26-37),
# fix the number format on E25, and update the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E25 changes from general to the 5-decimal number format (style 1 -> 2) ---
$ws.Range("E25").NumberFormat = "0.00000"

# --- New Panoche transect data for rows 26-37 ---
# Columns: A=Site, B=Site.Number, C=Site.Density, D=Transect number,
#          E=start.lat, F=start.lng, G=end.lat, H=end.lng
$rows = @(
    @{ R=26; A="Panoche"; B=1; C="High";   D=1; E=36.695819999999998;  F=-120.79666;            G=36.69623;             H=-120.79564999999999  },
    @{ R=27; A="Panoche"; B=1; C="High";   D=2; E=36.695500000000003;  F=-120.79646;            G=36.695999999999998;   H=-120.79554           },
    @{ R=28; A="Panoche"; B=1; C="High";   D=3; E=36.69538;            F=-120.79626;            G=36.695909999999998;   H=-120.79533000000001  },
    @{ R=29; A="Panoche"; B=1; C="Medium"; D=1; E=36.695929999999997;  F=-120.79772;            G=36.695279999999997;   H=-120.79848           },
    @{ R=30; A="Panoche"; B=1; C="Medium"; D=2; E=36.696159999999999;  F=-120.79794;            G=36.695610000000002;   H=-120.79882000000001  },
    @{ R=31; A="Panoche"; B=1; C="Medium"; D=3; E=36.696379999999998;  F=-120.79819999999999;   G=36.695790000000002;   H=-120.79904999999999  },
    @{ R=32; A="Panoche"; B=1; C="Low";    D=1; E=36.695239999999998;  F=-120.79682;            G=36.694920000000003;   H=-120.79575           },
    @{ R=33; A="Panoche"; B=1; C="Low";    D=2; E=36.695129999999999;  F=-120.79688;            G=36.694839999999999;   H=-120.79583           },
    @{ R=34; A="Panoche"; B=1; C="Low";    D=3; E=36.694960000000002;  F=-120.79707000000001;   G=36.69464;             H=-120.79601           },
    @{ R=35; A="Panoche"; B=1; C="None";   D=1; E=36.693579999999997;  F=-120.79268999999999;   G=36.692770000000003;   H=-120.79217           },
    @{ R=36; A="Panoche"; B=1; C="None";   D=2; E=36.692860000000003;  F=-120.79191;            G=36.696370000000002;   H=-120.79243           },
    @{ R=37; A="Panoche"; B=1; C="None";   D=3; E=36.693759999999997;  F=-120.79217;            G=36.692999999999998;   H=-120.79161000000001  }
)

# Cells that need the 5-decimal "0.00000" number format applied (style 2),
# everything else keeps the plain bordered style (style 1).
$fiveDecimalCells = @("E27", "G27", "F31", "G37")

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("A$r").Value = $row.A
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
}

foreach ($addr in $fiveDecimalCells) {
    $ws.Range($addr).NumberFormat = "0.00000"
}

# --- View state: scroll so row 5 is at the top, and select G37 ---
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G37").Select()
